$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (Column A dates as OLE automation date serials,
# matching the existing values already used in the sheet).
$newRows = @(
    @{ Row = 234; Date = 44308; B = 2; C = 26; D = 167.9261125104954 },
    @{ Row = 235; Date = 44309; B = 5; C = 24; D = 155.0087192404573 },
    @{ Row = 236; Date = 44310; B = 2; C = 23; D = 148.5500226054382 },
    @{ Row = 237; Date = 44311; B = 3; C = 22; D = 142.0913259704192 },
    @{ Row = 238; Date = 44312; B = 0; C = 17; D = 109.7978427953239 }
)

# Column A in the existing data uses a specific date style (border, bold,
# centered, custom date number format). Copy that formatting down onto each
# new date cell before writing its value.
$lastExistingDateCell = $ws.Cells.Item(233, 1)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    $lastExistingDateCell.Copy()
    $cellA = $ws.Cells.Item($rowIndex, 1)
    $cellA.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $cellA.Value = [DateTime]::FromOADate($r.Date)

    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
}
